# "major accuracy check update"
#
# 1. Library-kit name typo/accuracy fix: the shared string used by every
#    row's column G ("NEBNextPoly(A)E7490") is missing the kit's "Large"
#    suffix -> append "L" so it reads "NEBNextPoly(A)E7490L".
# 2. The "accuracy check" boolean in column I (rows 2-27) is re-entered as
#    an explicit =FALSE() formula instead of a bare boolean literal.
# 3. Column G is widened so the longer kit name is fully visible, and the
#    sheet's saved selection is moved onto the column that was just edited.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update column G, rows 2-27 (all rows share this one string).
$ws.Range("G2:G27").Value = "NEBNextPoly(A)E7490L"

# 2. Re-enter column I, rows 2-27 as an explicit FALSE() formula. Looping
#    cell-by-cell (rather than assigning .Formula on the whole G2:G27 range
#    at once) keeps each row an independent formula instead of Excel
#    collapsing the identical formulas into one shared-formula group, which
#    matches the per-row <f> markup the target sheet uses.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 9).Formula = "=FALSE()"
}

# 3. Widen column G to fit the updated text (~31.66 chars) and move the
#    active selection from I2:I27 onto G2:G27.
$ws.Columns.Item(7).ColumnWidth = 30.75
$ws.Range("G2:G27").Select()
